$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 8.307976
$ws.Range("H2").Value = 24.923928
$ws.Range("I2").Value = 0.463503063770797
$ws.Range("J2").Value = 0.4635030637707969
$ws.Range("M2").Value = 71.05094633333333
$ws.Range("N2").Value = 213.152839
$ws.Range("O2").Value = 0.8240565632932695
$ws.Range("P2").Value = 0.8240565632932696
$ws.Range("Q2").Value = 590.2895569146212
$ws.Range("R2").Value = 5312.606012231592
$ws.Range("S2").Value = 0.3819527418068641
$ws.Range("T2").Value = 0.3819527418068641
$ws.Range("G3").Value = 8.307976
$ws.Range("H3").Value = 24.923928
$ws.Range("I3").Value = 0.463503063770797
$ws.Range("J3").Value = 0.4635030637707969
$ws.Range("O3").Value = 0.1323102827659759
$ws.Range("P3").Value = 0.132310282765976
$ws.Range("Q3").Value = 94.77672003126933
$ws.Range("R3").Value = 852.9904802814241
$ws.Range("S3").Value = 0.06132622143041032
$ws.Range("T3").Value = 0.06132622143041033
$ws.Range("G4").Value = 8.307976
$ws.Range("H4").Value = 24.923928
$ws.Range("I4").Value = 0.463503063770797
$ws.Range("J4").Value = 0.4635030637707969
$ws.Range("M4").Value = 3.762092333333333
$ws.Range("N4").Value = 11.286277
$ws.Range("O4").Value = 0.04363315394075456
$ws.Range("P4").Value = 0.04363315394075455
$ws.Range("Q4").Value = 31.25537281511733
$ws.Range("R4").Value = 281.298355336056
$ws.Range("S4").Value = 0.02022410053352256
$ws.Range("T4").Value = 0.02022410053352255
$ws.Range("I5").Value = 0.1830631876040579
$ws.Range("J5").Value = 0.1830631876040579
$ws.Range("M5").Value = 71.05094633333333
$ws.Range("N5").Value = 213.152839
$ws.Range("O5").Value = 0.8240565632932695
$ws.Range("P5").Value = 0.8240565632932696
$ws.Range("Q5").Value = 233.1382386538302
$ws.Range("R5").Value = 2098.244147884472
$ws.Range("S5").Value = 0.150854421242511
$ws.Range("T5").Value = 0.150854421242511
$ws.Range("I6").Value = 0.1830631876040579
$ws.Range("J6").Value = 0.1830631876040579
$ws.Range("O6").Value = 0.1323102827659759
$ws.Range("P6").Value = 0.132310282765976
$ws.Range("S6").Value = 0.0242211421159338
$ws.Range("T6").Value = 0.02422114211593381
$ws.Range("I7").Value = 0.1830631876040579
$ws.Range("J7").Value = 0.1830631876040579
$ws.Range("M7").Value = 3.762092333333333
$ws.Range("N7").Value = 11.286277
$ws.Range("O7").Value = 0.04363315394075456
$ws.Range("P7").Value = 0.04363315394075455
$ws.Range("Q7").Value = 12.34448836376622
$ws.Range("R7").Value = 111.100395273896
$ws.Range("S7").Value = 0.007987624245613089
$ws.Range("T7").Value = 0.007987624245613089
$ws.Range("G8").Value = 5.1529
$ws.Range("H8").Value = 15.4587
$ws.Range("I8").Value = 0.2874809625478624
$ws.Range("J8").Value = 0.2874809625478624
$ws.Range("M8").Value = 71.05094633333333
$ws.Range("N8").Value = 213.152839
$ws.Range("O8").Value = 0.8240565632932695
$ws.Range("P8").Value = 0.8240565632932696
$ws.Range("Q8").Value = 366.1184213610333
$ws.Range("R8").Value = 3295.0657922493
$ws.Range("S8").Value = 0.2369005740094326
$ws.Range("T8").Value = 0.2369005740094326
$ws.Range("G9").Value = 5.1529
$ws.Range("H9").Value = 15.4587
$ws.Range("I9").Value = 0.2874809625478624
$ws.Range("J9").Value = 0.2874809625478624
$ws.Range("O9").Value = 0.1323102827659759
$ws.Range("P9").Value = 0.132310282765976
$ws.Range("Q9").Value = 58.78386753273333
$ws.Range("R9").Value = 529.0548077946
$ws.Range("S9").Value = 0.0380366874445426
$ws.Range("T9").Value = 0.03803668744454261
$ws.Range("G10").Value = 5.1529
$ws.Range("H10").Value = 15.4587
$ws.Range("I10").Value = 0.2874809625478624
$ws.Range("J10").Value = 0.2874809625478624
$ws.Range("M10").Value = 3.762092333333333
$ws.Range("N10").Value = 11.286277
$ws.Range("O10").Value = 0.04363315394075456
$ws.Range("P10").Value = 0.04363315394075455
$ws.Range("Q10").Value = 19.38568558443333
$ws.Range("R10").Value = 174.4711702599
$ws.Range("S10").Value = 0.01254370109388717
$ws.Range("T10").Value = 0.01254370109388717
$ws.Range("G11").Value = 1.182158666666667
$ws.Range("H11").Value = 3.546476
$ws.Range("I11").Value = 0.06595278607728289
$ws.Range("J11").Value = 0.06595278607728287
$ws.Range("M11").Value = 71.05094633333333
$ws.Range("N11").Value = 213.152839
$ws.Range("O11").Value = 0.8240565632932695
$ws.Range("P11").Value = 0.8240565632932696
$ws.Range("Q11").Value = 83.99349198281823
$ws.Range("R11").Value = 755.941427845364
$ws.Range("S11").Value = 0.05434882623446192
$ws.Range("T11").Value = 0.05434882623446192
$ws.Range("G12").Value = 1.182158666666667
$ws.Range("H12").Value = 3.546476
$ws.Range("I12").Value = 0.06595278607728289
$ws.Range("J12").Value = 0.06595278607728287
$ws.Range("O12").Value = 0.1323102827659759
$ws.Range("P12").Value = 0.132310282765976
$ws.Range("Q12").Value = 13.48597070853423
$ws.Range("R12").Value = 121.373736376808
$ws.Range("S12").Value = 0.008726231775089219
$ws.Range("T12").Value = 0.008726231775089219
$ws.Range("G13").Value = 1.182158666666667
$ws.Range("H13").Value = 3.546476
$ws.Range("I13").Value = 0.06595278607728289
$ws.Range("J13").Value = 0.06595278607728287
$ws.Range("M13").Value = 3.762092333333333
$ws.Range("N13").Value = 11.286277
$ws.Range("O13").Value = 0.04363315394075456
$ws.Range("P13").Value = 0.04363315394075455
$ws.Range("Q13").Value = 4.447390056650223
$ws.Range("R13").Value = 40.026510509852
$ws.Range("S13").Value = 0.002877728067731738
$ws.Range("T13").Value = 0.002877728067731737
